# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 14:22"

# Paises Bajos (row 17)
$ws.Range("B17").Value = 37190
$ws.Range("C17").Value = 655
$ws.Range("E17").Value = 32531
$ws.Range("G17").Value = 120
$ws.Range("H17").Value = 4409

# Dinamarca (row 40)
$ws.Range("B40").Value = 8445
$ws.Range("C40").Value = 235
$ws.Range("D40").Value = 5669
$ws.Range("E40").Value = 2358
$ws.Range("F40").Value = 70
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 418

# Croacia (row 66)
$ws.Range("B66").Value = 2016
$ws.Range("C66").Value = 7
$ws.Range("D66").Value = 1034
$ws.Range("E66").Value = 928
$ws.Range("F66").Value = 21
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 54

# Azerbaiyan (row 73)
$ws.Range("B73").Value = 1617
$ws.Range("C73").Value = 25
$ws.Range("D73").Value = 1080
$ws.Range("E73").Value = 516
$ws.Range("F73").Value = 15

# Kenia (row 116)
$ws.Range("B116").Value = 343
$ws.Range("C116").Value = 7
$ws.Range("E116").Value = 235

# Vietnam (row 126)
$ws.Range("D126").Value = 225
$ws.Range("E126").Value = 45

# Togo / Cabo Verde swap places (row 149 / 150) with updated Cabo Verde numbers
$ws.Range("A149").Value = "Cabo Verde"
$ws.Range("B149").Value = 90
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 88
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

$ws.Range("A150").Value = "Togo"
$ws.Range("B150").Value = 90
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 59
$ws.Range("E150").Value = 25
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6
